$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 518.58826
$ws.Range("I107").Value = 457.85715
$ws.Range("J107").Value = 802
$ws.Range("K107").Value = 457.85715
$ws.Range("L107").Value = 802
$ws.Range("M107").Value = 1462.14285
$ws.Range("N107").Value = -4642
$ws.Range("H113").Value = 4672.727
$ws.Range("I113").Value = 2900
$ws.Range("J113").Value = 5066.6665
$ws.Range("K113").Value = 2900
$ws.Range("L113").Value = 5066.6665
$ws.Range("M113").Value = 354
$ws.Range("N113").Value = -11574.6665
$ws.Range("H123").Value = 29999.6
$ws.Range("J123").Value = 29999.6
$ws.Range("L123").Value = 29999.6
$ws.Range("N123").Value = -39799.6
$ws.Range("H132").Value = 2775.0356
$ws.Range("I132").Value = 1920.875
$ws.Range("J132").Value = 7900
$ws.Range("K132").Value = 5762.625
$ws.Range("L132").Value = 23700
$ws.Range("M132").Value = -3232.625
$ws.Range("N132").Value = -28760
$ws.Range("H138").Value = 6062670.5
$ws.Range("I138").Value = 1452.84
$ws.Range("J138").Value = 25003976
$ws.Range("K138").Value = 4358.52
$ws.Range("L138").Value = 75011928
$ws.Range("M138").Value = 781.4800000000005
$ws.Range("N138").Value = -75022208

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 112.5
$ws.Range("I5").Value = 112.5
$ws.Range("K5").Value = 112.5
$ws.Range("M5").Value = -0.5
$ws.Range("H44").Value = 30125
$ws.Range("J44").Value = 30125
$ws.Range("L44").Value = 30125
$ws.Range("N44").Value = -31101
$ws.Range("H55").Value = 22969.6
$ws.Range("I55").Value = 15424
$ws.Range("J55").Value = 28000
$ws.Range("K55").Value = 15424
$ws.Range("L55").Value = 28000
$ws.Range("M55").Value = -15109
$ws.Range("N55").Value = -28630

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 112.5
$ws.Range("I4").Value = 112.5
$ws.Range("K4").Value = 112.5
$ws.Range("M4").Value = 2.5
$ws.Range("H82").Value = 10263.857
$ws.Range("I82").Value = 4369.4
$ws.Range("J82").Value = 25000
$ws.Range("K82").Value = 4369.4
$ws.Range("L82").Value = 25000
$ws.Range("M82").Value = -3986.4
$ws.Range("N82").Value = -25766
$ws.Range("H85").Value = 10263.857
$ws.Range("I85").Value = 4369.4
$ws.Range("J85").Value = 25000
$ws.Range("K85").Value = 4369.4
$ws.Range("L85").Value = 25000
$ws.Range("M85").Value = -3043.4
$ws.Range("N85").Value = -27652

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 928.1111
$ws.Range("I16").Value = 835
$ws.Range("J16").Value = 1002.6
$ws.Range("K16").Value = 835
$ws.Range("L16").Value = 1002.6
$ws.Range("M16").Value = -548
$ws.Range("N16").Value = -1576.6
$ws.Range("H31").Value = 60001532
$ws.Range("I31").Value = 71429736
$ws.Range("J31").Value = 45456544
$ws.Range("K31").Value = 71429736
$ws.Range("L31").Value = 45456544
$ws.Range("M31").Value = -71429441
$ws.Range("N31").Value = -45457134
$ws.Range("H34").Value = 60001532
$ws.Range("I34").Value = 71429736
$ws.Range("J34").Value = 45456544
$ws.Range("K34").Value = 71429736
$ws.Range("L34").Value = 45456544
$ws.Range("M34").Value = -71429534
$ws.Range("N34").Value = -45456948
$ws.Range("H50").Value = 12499.5
$ws.Range("I50").Value = 5000
$ws.Range("J50").Value = 14999.333
$ws.Range("K50").Value = 5000
$ws.Range("L50").Value = 14999.333
$ws.Range("M50").Value = -4375
$ws.Range("N50").Value = -16249.333
$ws.Range("H51").Value = 17998
$ws.Range("J51").Value = 17998
$ws.Range("L51").Value = 17998
$ws.Range("N51").Value = -19470
$ws.Range("H58").Value = 6972.737
$ws.Range("I58").Value = 9881.166999999999
$ws.Range("J58").Value = 1986.8572
$ws.Range("K58").Value = 9881.166999999999
$ws.Range("L58").Value = 1986.8572
$ws.Range("M58").Value = -9678.166999999999
$ws.Range("N58").Value = -2392.8572
$ws.Range("H59").Value = 25000
$ws.Range("J59").Value = 25000
$ws.Range("L59").Value = 25000
$ws.Range("N59").Value = -27290
$ws.Range("H61").Value = 17998
$ws.Range("J61").Value = 17998
$ws.Range("L61").Value = 17998
$ws.Range("N61").Value = -18694
$ws.Range("H68").Value = 22497.5
$ws.Range("J68").Value = 22497.5
$ws.Range("L68").Value = 22497.5
$ws.Range("N68").Value = -23995.5
$ws.Range("H71").Value = 22497.5
$ws.Range("J71").Value = 22497.5
$ws.Range("L71").Value = 67492.5
$ws.Range("N71").Value = -74980.5
$ws.Range("H74").Value = 26292.309
$ws.Range("J74").Value = 26292.309
$ws.Range("L74").Value = 26292.309
$ws.Range("N74").Value = -28040.309
$ws.Range("H77").Value = 26292.309
$ws.Range("J77").Value = 26292.309
$ws.Range("L77").Value = 78876.927
$ws.Range("N77").Value = -87612.927
$ws.Range("H113").Value = 928.1111
$ws.Range("I113").Value = 835
$ws.Range("J113").Value = 1002.6
$ws.Range("K113").Value = 835
$ws.Range("L113").Value = 1002.6
$ws.Range("M113").Value = 1335
$ws.Range("N113").Value = -5342.6
$ws.Range("H132").Value = 2477.4688
$ws.Range("I132").Value = 2167.9614
$ws.Range("K132").Value = 6503.8842
$ws.Range("M132").Value = -3973.8842
$ws.Range("H134").Value = 4222.4
$ws.Range("I134").Value = 4222.4
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 12667.2
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -10132.2
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 6972.737
$ws.Range("I136").Value = 9881.166999999999
$ws.Range("J136").Value = 1986.8572
$ws.Range("K136").Value = 29643.501
$ws.Range("L136").Value = 5960.571599999999
$ws.Range("M136").Value = -27093.501
$ws.Range("N136").Value = -11060.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 713.7
$ws.Range("I5").Value = 591
$ws.Range("K5").Value = 1773
$ws.Range("M5").Value = -1661
$ws.Range("H100").Value = 3866.6667
$ws.Range("J100").Value = 3866.6667
$ws.Range("L100").Value = 11600.0001
$ws.Range("N100").Value = -13222.0001
$ws.Range("H131").Value = 1107.1666
$ws.Range("I131").Value = 840.5
$ws.Range("J131").Value = 1131.409
$ws.Range("K131").Value = 2521.5
$ws.Range("L131").Value = 3394.227
$ws.Range("M131").Value = 2518.5
$ws.Range("N131").Value = -13474.227
$ws.Range("H135").Value = 713.7
$ws.Range("I135").Value = 591
$ws.Range("K135").Value = 5319
$ws.Range("M135").Value = -2784

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4081.5454
$ws.Range("I132").Value = 4269.8
$ws.Range("J132").Value = 3924.6667
$ws.Range("K132").Value = 12809.4
$ws.Range("L132").Value = 11774.0001
$ws.Range("M132").Value = -10279.4
$ws.Range("N132").Value = -16834.0001
$ws.Range("H134").Value = 20000
$ws.Range("J134").Value = 20000
$ws.Range("L134").Value = 60000
$ws.Range("N134").Value = -65070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 77737.08
$ws.Range("I22").Value = 333826.66
$ws.Range("J22").Value = 910.2
$ws.Range("K22").Value = 333826.66
$ws.Range("L22").Value = 910.2
$ws.Range("M22").Value = -333531.66
$ws.Range("N22").Value = -1500.2
$ws.Range("H27").Value = 77737.08
$ws.Range("I27").Value = 333826.66
$ws.Range("J27").Value = 910.2
$ws.Range("K27").Value = 333826.66
$ws.Range("L27").Value = 910.2
$ws.Range("M27").Value = -333719.66
$ws.Range("N27").Value = -1124.2
$ws.Range("H141").Value = 40695
$ws.Range("J141").Value = 40695
$ws.Range("L141").Value = 40695
$ws.Range("M141").Value = -51055

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 13499.5
$ws.Range("I51").Value = 9999
$ws.Range("J51").Value = 17000
$ws.Range("K51").Value = 9999
$ws.Range("L51").Value = 17000
$ws.Range("M51").Value = -9489
$ws.Range("N51").Value = -18020
$ws.Range("H57").Value = 24285.762
$ws.Range("J57").Value = 24285.762
$ws.Range("L57").Value = 24285.762
$ws.Range("N57").Value = -25793.762
$ws.Range("H100").Value = 928.65515
$ws.Range("I100").Value = 943.24
$ws.Range("J100").Value = 837.5
$ws.Range("K100").Value = 1886.48
$ws.Range("L100").Value = 1675
$ws.Range("M100").Value = -1345.48
$ws.Range("N100").Value = -2757
$ws.Range("H124").Value = 30976.334
$ws.Range("J124").Value = 30976.334
$ws.Range("L124").Value = 30976.334
$ws.Range("N124").Value = -40796.334
